$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.505614041169197
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 3.082599426703578
$ws.Range("E2").Value = 71517.89157740913
$ws.Range("G2").Value = 71524.13302733591

# Row 3
$ws.Range("B3").Value = 0.7287194209349384
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 3.594575437922795

# Row 4
$ws.Range("B4").Value = 0.00006486019690155054
$ws.Range("C4").Value = 0.004309184025731883
$ws.Range("D4").Value = 16.98373111632243
$ws.Range("E4").Value = 198602002.3250627
$ws.Range("G4").Value = 198602019.3131678
